# Apply updated dSF (column F) values as part of a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -2
$ws.Range("F3").Value  = -4
$ws.Range("F4").Value  = -3
$ws.Range("F5").Value  = 4
$ws.Range("F6").Value  = 3
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -3
